$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 349  # was 347
$ws.Range("F6").Value = 405  # was 403
$ws.Range("F7").Value = 897  # was 896
$ws.Range("F9").Value = 544  # was 539
$ws.Range("F12").Value = 1172  # was 1169
$ws.Range("F14").Value = 254  # was 252
$ws.Range("F15").Value = 43  # was 42
$ws.Range("F17").Value = 6726  # was 6716
$ws.Range("F18").Value = 66  # was 65
$ws.Range("F21").Value = 7643  # was 7635
$ws.Range("F23").Value = 38  # was 37
$ws.Range("F24").Value = 3419  # was 3416
$ws.Range("F26").Value = 2151  # was 2145
$ws.Range("F27").Value = 917  # was 916
$ws.Range("F28").Value = 4525  # was 4522
$ws.Range("F29").Value = 184  # was 173
$ws.Range("F30").Value = 354  # was 353
$ws.Range("F32").Value = 2  # was 1
$ws.Range("F33").Value = 241  # was 240
$ws.Range("F35").Value = 1775  # was 1767
$ws.Range("F37").Value = 200  # was 195
$ws.Range("F39").Value = 6  # was 5
$ws.Range("F41").Value = 1246  # was 1243
$ws.Range("F42").Value = 1875  # was 1861
$ws.Range("F43").Value = 2151  # was 2150

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 72  # was 71
$ws.Range("F4").Value = 55  # was 54

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 244  # was 243
$ws.Range("F3").Value = 1241  # was 1239
$ws.Range("F4").Value = 79  # was 78

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 244  # was 243
$ws.Range("F4").Value = 1241  # was 1239
$ws.Range("F5").Value = 79  # was 78
$ws.Range("F7").Value = 349  # was 347
$ws.Range("F8").Value = 405  # was 403
$ws.Range("F9").Value = 897  # was 896
$ws.Range("F11").Value = 544  # was 539
$ws.Range("F14").Value = 1172  # was 1169
$ws.Range("F15").Value = 72  # was 71
$ws.Range("F17").Value = 254  # was 252
$ws.Range("F18").Value = 43  # was 42
$ws.Range("F20").Value = 6726  # was 6716
$ws.Range("F21").Value = 66  # was 65
$ws.Range("F24").Value = 7643  # was 7635
$ws.Range("F26").Value = 38  # was 37
$ws.Range("F27").Value = 3419  # was 3416
$ws.Range("F29").Value = 2151  # was 2145
$ws.Range("F30").Value = 917  # was 916
$ws.Range("F31").Value = 4525  # was 4522
$ws.Range("F32").Value = 184  # was 173
$ws.Range("F33").Value = 354  # was 353
$ws.Range("F35").Value = 55  # was 54
$ws.Range("F36").Value = 2  # was 1
$ws.Range("F37").Value = 241  # was 240
$ws.Range("F38").Value = 1775  # was 1767
$ws.Range("F40").Value = 200  # was 195
$ws.Range("F42").Value = 6  # was 5
$ws.Range("F44").Value = 1246  # was 1243
$ws.Range("F45").Value = 1875  # was 1861
$ws.Range("F47").Value = 2151  # was 2150
